$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.320.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.071.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +6.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.69%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +5.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.61'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.88%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.113'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.371'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.94%  '
$ws.Range("E12").Value = '  +2.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.598.38'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.45'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000174'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +16.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '58.313.60'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +9.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.074.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '342.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.97%  '
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.507'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.56%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0₃0977'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.63%  '
$ws.Range("E27").Value = '  +4.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("E29").Value = '  +9.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.54'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.26%  '
$ws.Range("E31").Value = '  +7.44%  '
$ws.Range("E32").Value = '  +6.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  +9.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = '  +3.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '26.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +15.06%  '
$ws.Range("E39").Value = '  +4.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.108.02'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.97'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +11.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.670'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.89%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.84%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.341.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.82%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.07%  '
$ws.Range("E48").Value = '  +4.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.99%  '
$ws.Range("E50").Value = '  +3.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.56%  '
